$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 value from 9 to 18
$ws.Range("B1").Value = 18.0

# Row 12
$ws.Range("A12").Value = "Прямоугольник"
$ws.Range("B12").Value = 50.0
$ws.Range("C12").Value = 50.0

# Row 13
$ws.Range("A13").Value = "Отрезок"
$ws.Range("B13").Value = 100.0
$ws.Range("C13").Value = 100.0
$ws.Range("D13").Value = 125.0
$ws.Range("E13").Value = 125.0

# Row 14
$ws.Range("A14").Value = "Прямоугольник"
$ws.Range("B14").Value = 125.0
$ws.Range("C14").Value = 125.0

# Row 15
$ws.Range("A15").Value = "Прямоугольник"
$ws.Range("B15").Value = 200.0
$ws.Range("C15").Value = 50.0

# Row 16
$ws.Range("A16").Value = "Прямоугольник"
$ws.Range("B16").Value = 50.0
$ws.Range("C16").Value = 200.0

# Row 17
$ws.Range("A17").Value = "Прямоугольник"
$ws.Range("B17").Value = 200.0
$ws.Range("C17").Value = 200.0

# Row 18
$ws.Range("A18").Value = "Отрезок"
$ws.Range("B18").Value = 200.0
$ws.Range("C18").Value = 100.0
$ws.Range("D18").Value = 175.0
$ws.Range("E18").Value = 125.0

# Row 19
$ws.Range("A19").Value = "Отрезок"
$ws.Range("B19").Value = 100.0
$ws.Range("C19").Value = 200.0
$ws.Range("D19").Value = 125.0
$ws.Range("E19").Value = 175.0

# Row 20
$ws.Range("A20").Value = "Отрезок"
$ws.Range("B20").Value = 200.0
$ws.Range("C20").Value = 200.0
$ws.Range("D20").Value = 175.0
$ws.Range("E20").Value = 175.0
